$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.493.14'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.81%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.601.64'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.53%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '537.88'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.93%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.41'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.37%  '

$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.51'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.67%  '

$ws.Range('E10').Value = '  +1.53%  '

$ws.Range('E11').Value = '  +1.72%  '

$ws.Range('E12').Value = '  -0.97%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.058.99'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.38%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '59.395.76'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.73%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.77'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.15%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.595.33'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.59%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000133'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.36%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '341.50'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.01%  '

$ws.Range('E19').Value = '  +1.46%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.10'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.12%  '

$ws.Range('E21').Value = '  -1.97%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.02%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.41'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.72%  '

$ws.Range('E24').Value = '  +1.22%  '

$ws.Range('E25').Value = '  -1.75%  '

$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('E27').Value = '  +2.72%  '

$ws.Range('E28').Value = '  +2.67%  '

$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('E30').Value = '  +5.33%  '

$ws.Range('E31').Value = '  -1.35%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.80'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.51%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '149.88'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.54%  '

$ws.Range('E34').Value = '  -0.16%  '

$ws.Range('E35').Value = '  -0.47%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.840'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.78%  '

$ws.Range('E37').Value = '  -0.54%  '

$ws.Range('E38').Value = '  -0.02%  '

$ws.Range('E39').Value = '  +0.40%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '271.62'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.19%  '

$ws.Range('E42').Value = '  +1.67%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '10.76'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.12%  '

$ws.Range('E44').Value = '  -0.09%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0525'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.63%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '18.62'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.73%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0223'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.43%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.941.80'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.08%  '

$ws.Range('E49').Value = '  -0.42%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '112.02'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.09%  '

$ws.Range('E51').Value = '  +1.09%  '

